# Armedospor vs Of 1461 - Maç linki eklendi.
# Ajans Of vs Kural Kesiciler maç sonucu eklendi.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Maç linki eklendi: turn the existing video URL text in H12
# (Narin Güran Grubu / Of FK - 61.Alay) into a real clickable hyperlink,
# and add the new match-video URL as text for H13 (Armedospor - Of 1461).
$ws.Hyperlinks.Add($ws.Range("H12"), "https://youtu.be/j0ONT2EiueM")
$ws.Range("H13").Value = "https://youtu.be/r6AEVdp_RAM"

# --- Maç sonucu eklendi: Ajans Of 3 - 1 Kural Kesiciler (row 14)
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 1

# Leave the selection where the author's edits ended.
$ws.Range("H16").Select()
